$d = $word.ActiveDocument

# Change 1: merge " #Region" + "Table" + "_T" runs into " #Region_T"
$d.Content.Find.Execute("#RegionTable_T", $true, $false, $false, $false, $false,
                         $true, 1, $false, "#Region_T", 2)

# Change 2: remove proofErr markers and merge runs so that text reads
# "...Lorem ipsum dolor sit amet. Lorem ipsum..."
$d.Content.Find.Execute("dolor sit amet.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dolor sit amet.", 2)
